# Surat Keterangan Usaha Update: rename JULAEHA -> KIRMAN, update NIK,
# postal code, pekerjaan, and the letter date.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2) | Out-Null
}

# Name (occurs in both the body paragraph and the signature block table)
Replace-Text "JULAEHA" "KIRMAN"

# NIK
Replace-Text "3208074101710002" "3208270107850249"

# Postal code in "Tempat, Tgl Lahir"
Replace-Text "KUNINGAN, 25934" "KUNINGAN, 31229"

# Pekerjaan / kegiatan usaha
Replace-Text "Mengurus Rumah Tangga" "m"

# Letter date in the signature block
Replace-Text "Ciawigebang, 02 Oktober 2017" "Ciawigebang, 16 Oktober 2017"
